$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsReliability = $wb.Worksheets.Item("ESUfRaLCD-reliability")

# ---------------------------------------------------------------------
# ESUfRaLCD-reliability: insert 4 new rows at row 8 for the "w ccs"
# sources (hard coal w ccs, natural gas combined cycle w ccs,
# biomass w ccs, lignite w ccs), pushing the hydrogen rows down.
# ---------------------------------------------------------------------
$wsReliability.Rows("8:11").Insert()

# Restore the "power plants" formula for the 4 freshly inserted rows
# (row insert does not auto-fill the formula into brand-new rows).
$wsReliability.Range("C8:C11").Formula = "=IF(A8=`"`",`"`",CONCATENATE(A8,`" power plants`"))"

# Column B literal ("es") labels are typed first for the two new plant
# types so their shared-string entries are created in the same order
# as the source edit, then the base-name column A, then the two rows
# that reuse the already-existing "w CCS es" strings.
$wsReliability.Range("B8").Value = "hard coal w ccs es"
$wsReliability.Range("B9").Value = "natural gas combined cycle w ccs es"

$wsReliability.Range("A8").Value = "hard coal w ccs"
$wsReliability.Range("A9").Value = "natural gas combined cycle w ccs"
$wsReliability.Range("A10").Value = "biomass w ccs"
$wsReliability.Range("A11").Value = "lignite w ccs"

$wsReliability.Range("B10").Value = "biomass w CCS es"
$wsReliability.Range("B11").Value = "lignite w CCS es"

# ---------------------------------------------------------------------
# View / selection state: reliability sheet becomes the active tab,
# About sheet loses tabSelected and gets a new resting selection.
# ---------------------------------------------------------------------
$wsAbout.Range("A29:A30").Select()

$wsReliability.Activate()
$wsReliability.Range("A12").Select()
